$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Hydrogen row (row 3): update B3, clear D3 (value removed, cell stays present but empty)
$ws.Range("B3").Value2 = 1907824.99868752
$ws.Range("D3").ClearContents()
$ws.Range("D3").Style = "Normal"

# 2. Methanol row (row 4): update C4
$ws.Range("C4").Value2 = 48.91406560215935

# 3. Ammonia row (row 5): update C5
$ws.Range("C5").Value2 = 3071.03425598076

# 4. Row 7: rename "Other" -> "Biogas", update D7
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value2 = 1705.024607732538

# 5. New row 8: "Other" with D8 value, formatted like row 7's A/D cells
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value2 = 1285.009340266951
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Style = "Normal"
